$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new rows (16 and 17) with the same look & feel as the existing
#     data rows. Copy formats from row 15 (A:E) down into rows 16 and 17 so
#     that column A keeps its bordered/bold/centered style (s="1") and the
#     other columns keep the default (unstyled) look.
$ws.Range("A15:E15").Copy()
$ws.Range("A16:E17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update row 8 (index 6) ---
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true
$ws.Range("B8").Value = "line7"

# --- Update row 9 (index 7) ---
$ws.Range("C9").Value = 16
$ws.Range("B9").Value = "line8"

# --- Update row 10 (index 8) ---
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("B10").Value = "extr1"

# --- Update row 11 (index 9) ---
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("B11").Value = "extr2"

# --- Update row 12 (index 10) ---
$ws.Range("C12").Value = 10
$ws.Range("E12").Value = $true
$ws.Range("B12").Value = "extr3"

# --- Update row 13 (index 11) ---
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true
$ws.Range("B13").Value = "extr4"

# --- Update row 14 (index 12) ---
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false
$ws.Range("B14").Value = "extr5"

# --- Update row 15 (index 13) ---
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $true
$ws.Range("B15").Value = "extr6"

# --- New row 16 (index 14) ---
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

# --- New row 17 (index 15) ---
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false
